$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pl_mw values per row (row number -> column letter -> value)
$newValues = @{
    2 = @{ "B"=0.09261994809251917; "D"=0.01438135605270219; "E"=0.4212590745016769; "F"=0.5722255841690824; "G"=0.00239282271308533; "I"=0.3798438538984508; "K"=0.7332698380719194; "N"=1.111352089838576; "O"=1.870665470003445 }
    3 = @{ "B"=0.08232607355142818; "D"=0.01317700204472771; "E"=0.3674750499019979; "F"=0.5633032175208612; "G"=0.002395681757626137; "I"=0.3856816045368272; "K"=0.6393757388666472; "N"=1.122742794153986; "O"=1.85481449600428 }
    4 = @{ "B"=0.07600601176710597; "D"=0.01243209405505397; "E"=0.3345423279487392; "F"=0.5582619178721586; "G"=0.002397529485395561; "I"=0.3894717625190145; "K"=0.5815714157743059; "N"=1.130171987566406; "O"=1.846519801055393 }
    5 = @{ "B"=0.0734308672458468; "D"=0.01212719614532887; "E"=0.3211430391357766; "F"=0.5563171463126295; "G"=0.00239830572416374; "I"=0.3910680011108543; "K"=0.5579778327977465; "N"=1.133308931003171; "O"=1.843500082263773 }
    6 = @{ "B"=0.07300329318718468; "D"=0.01207648772929204; "E"=0.3189193174280405; "F"=0.5560008290847733; "G"=0.002398436026113676; "I"=0.391336178382987; "K"=0.554057868385911; "N"=1.133836431998432; "O"=1.843020396929546 }
    7 = @{ "B"=0.0759712808224009; "D"=0.012427987496757; "E"=0.3343615375199533; "F"=0.5582352466480529; "G"=0.002397539859554762; "I"=0.3894930805723797; "K"=0.5812533769063464; "N"=1.130213850125049; "O"=1.846477618179847 }
    8 = @{ "B"=0.08907066583533663; "D"=0.01396723543883382; "E"=0.4026943835490897; "F"=0.5690582494689735; "G"=0.002393789407672495; "I"=0.3818139941085956; "K"=0.7009272389040007; "N"=1.115189290558533; "O"=1.864900934478783 }
    9 = @{ "B"=0.1147531366853087; "D"=0.01694172030767049; "E"=0.5375031219361119; "F"=0.5937659089395453; "G"=0.002387163428471938; "I"=0.368388604760602; "K"=0.934378178518898; "N"=1.089177883342366; "O"=1.91249366916972 }
    10 = @{ "B"=0.1336089052774838; "D"=0.01909927073328532; "E"=0.6371694737967175; "F"=0.6140673000649599; "G"=0.002382734727222127; "I"=0.3595213053098387; "K"=1.105141558108187; "N"=1.072168846167614; "O"=1.954533754587516 }
    11 = @{ "B"=0.1421820543462644; "D"=0.02007456286612097; "E"=0.6826740841936072; "F"=0.6237750530323041; "G"=0.002380814386658119; "I"=0.3557038172615288; "K"=1.182663082576482; "N"=1.064886816015175; "O"=1.975213628725442 }
    12 = @{ "B"=0.1454276391707481; "D"=0.02044297004265161; "E"=0.6999315022958115; "F"=0.6275194718587755; "G"=0.002380100686111709; "I"=0.3542893586653708; "K"=1.211995163569384; "N"=1.06219477921514; "O"=1.9832696502381 }
    13 = @{ "B"=0.1447286874371514; "D"=0.02036366804561851; "E"=0.6962136292964232; "F"=0.6267100008326025; "G"=0.002380253795187989; "I"=0.3545926019933594; "K"=1.205679037300342; "N"=1.062771643660284; "O"=1.981524613974102 }
    14 = @{ "B"=0.1424490896027493; "D"=0.02010489042175578; "E"=0.6840933338302051; "F"=0.6240817376988161; "G"=0.002380755400081433; "I"=0.3555868245603628; "K"=1.185076728191063; "N"=1.064664027513373; "O"=1.97587188558748 }
    15 = @{ "B"=0.141052648218178; "D"=0.01994626178702674; "E"=0.6766727224504763; "F"=0.6224807562351771; "G"=0.002381064402490999; "I"=0.3561998712621302; "K"=1.172454120313375; "N"=1.065831699093287; "O"=1.972438766990081 }
    16 = @{ "B"=0.1330485198012212; "D"=0.01903540643929347; "E"=0.6341991452694344; "F"=0.6134424144683948; "G"=0.002382862117705184; "I"=0.3597751427787568; "K"=1.100072073727972; "N"=1.072653910258872; "O"=1.953213690720929 }
    17 = @{ "B"=0.1281369373583487; "D"=0.01847502340999085; "E"=0.6081868541658082; "F"=0.6080189650599408; "G"=0.002383989060109129; "I"=0.362023880959657; "K"=1.055626664830129; "N"=1.076955780748861; "O"=1.941819088464825 }
    18 = @{ "B"=0.1253115221152257; "D"=0.0181521249100598; "E"=0.5932407683809799; "F"=0.60494399773971; "G"=0.002384646127612179; "I"=0.3633376516695801; "K"=1.030047860227398; "N"=1.079472972380515; "O"=1.935411559184644 }
    19 = @{ "B"=0.1243548232724407; "D"=0.01804269800016556; "E"=0.5881828827227622; "F"=0.6039104915535773; "G"=0.002384870126640348; "I"=0.3637859668377486; "K"=1.021384764591915; "N"=1.080332612720838; "O"=1.933267175091771 }
    20 = @{ "B"=0.1286598271580175; "D"=0.0185347374528817; "E"=0.6109542890941668; "F"=0.6085916967970206; "G"=0.002383868176622313; "I"=0.3617823918915706; "K"=1.060359509965508; "N"=1.076493402361592; "O"=1.94301690777769 }
    21 = @{ "B"=0.143118688578781; "D"=0.02018092470712674; "E"=0.6876526422603462; "F"=0.6248518657428832; "G"=0.002380607701177496; "I"=0.3552939518155735; "K"=1.191128772762568; "N"=1.064106410446975; "O"=1.977526112010167 }
    22 = @{ "B"=0.1525631448808298; "D"=0.02125145973987941; "E"=0.7379308102542836; "F"=0.6358771155677516; "G"=0.002378555393094376; "I"=0.3512349145535651; "K"=1.276455951213904; "N"=1.056392646890174; "O"=2.00139196711865 }
    23 = @{ "B"=0.1475230208738054; "D"=0.02068059215532259; "E"=0.7110819050139412; "F"=0.6299561720463345; "G"=0.00237964357928398; "I"=0.3533846759003225; "K"=1.230928091561509; "N"=1.060474682437906; "O"=1.988533821081575 }
    24 = @{ "B"=0.1284234338982628; "D"=0.01850774299587954; "E"=0.6097031047663251; "F"=0.6083326307149406; "G"=0.002383922799450877; "I"=0.361891503870817; "K"=1.058219873326493; "N"=1.07670230660203; "O"=1.942474927125289 }
    25 = @{ "B"=0.107806946252623; "D"=0.01614185482583963; "E"=0.5009339214195876; "F"=0.5867063062207336; "G"=0.00238887842224389; "I"=0.3718456092546853; "K"=0.8713555094118703; "N"=1.095845424461217; "O"=1.89838249409928 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Range("$col$row").Value = $newValues[$row][$col]
    }
}
